# Applies the "Pass xml instruction via XML instead of tag metadata" edit:
#  1. Remove the ":remove-old-xml" tag-metadata run from the [[Html2]] line
#     and replace it with the real "_GoBack" bookmark (which Word relocates
#     here automatically, vacating its previous position).
#  2. Merge the now-bookmark-free "HTML can also be imported..." sentence
#     back into a single run.
#  3. Merge the split "[[Html" / "3" / "]]" runs into a single "[[Html3]]" run.

$d = $word.ActiveDocument

# 1) Drop ":remove-old-xml" and drop a "_GoBack" bookmark in its place.
$r = $d.Content
$r.Find.Execute(":remove-old-xml", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0)
$r.Delete()
$d.Bookmarks.Add("_GoBack", $r)

# 2) Re-merge "...must be used" + ": " into one run (self-replace across the
#    run boundary coalesces the text into a single run).
$d.Content.Find.Execute("used: ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "used: ", 2)

# 3) Re-merge "[[Html" + "3" + "]]" into a single "[[Html3]]" run.
$d.Content.Find.Execute("[[Html3]]", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[[Html3]]", 2)
